$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters in order (matching the data columns below), skipping G (constant 1) and H (unused)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Each string: "<row>|<valB>|<valC>|...|<valN>"
$data = @(
    "2|1.02|1.035075658284013|1.042198635323328|1.038721552936539|1.051160215472998|1.038066846770822|1.04019141425881|1.044975999440203|1.041508779020077|1.053912488201351|1.04166860537868",
    "3|1.019999999999999|1.036142366995566|1.043007073317483|1.039736713705939|1.052110734089185|1.038300004548537|1.040900928344778|1.045595226354232|1.042333469880171|1.054675242628761|1.042379127056033",
    "4|1.02|1.036832733972142|1.043529991498238|1.040394077703443|1.052725888781589|1.038449147375976|1.041359605787217|1.045995061589079|1.042866980782597|1.055168278512431|1.042838455873051",
    "5|1.02|1.037122996708279|1.043749778848836|1.040670549762062|1.052984524439155|1.03851143348078|1.041552331582085|1.046162949309768|1.043091239886504|1.055375426418294|1.043031455360653",
    "6|1.02|1.037171735005649|1.043786679308836|1.040716977442259|1.053027951929076|1.038521867330893|1.041584685088901|1.046191126478026|1.043128892273552|1.055410200153372|1.04306385481316",
    "7|1.02|1.036836612349342|1.043532928494929|1.040397771480379|1.052729344587974|1.038449981271356|1.041362181400699|1.045997305711333|1.042869977457622|1.055171046922319|1.042841035144199",
    "8|1.02|1.035436129579981|1.042471890735106|1.039064530188107|1.051481425435626|1.03814600066041|1.040431285807705|1.045185444953681|1.041787511643725|1.054170371397386|1.041908817572689",
    "9|1.02|1.032969329180671|1.040600741845917|1.036718935441884|1.049283269394239|1.037597146020871|1.038787682291546|1.043748384094026|1.039879169741969|1.052403115723821|1.040262879951013",
    "10|1.02|1.031325476508602|1.039352363556767|1.035157743790574|1.047818430327238|1.037222386920933|1.037689779314523|1.042786027135653|1.038606354348663|1.051222327747684|1.039163417825831",
    "11|1.02|1.03061382820762|1.03881158521866|1.034482335737337|1.047184287437364|1.037058015417144|1.037213863686878|1.042368296958741|1.038055074660174|1.050710418226707|1.038686826343407",
    "12|1.02|1.030349512578942|1.038610683192909|1.034231549145788|1.046948760326361|1.036996645514211|1.037037009964147|1.042212979938539|1.037850283598582|1.050520179364937|1.038509721468104",
    "13|1.02|1.030406208178909|1.038653778838591|1.034285339644401|1.046999280697873|1.037009823817588|1.037074949191823|1.042246302904988|1.037894212935573|1.050560990443596|1.038547714573835",
    "14|1.02|1.030591979344067|1.038794979262114|1.034461603801937|1.047164818245752|1.037052948986483|1.03719924648267|1.042355461533662|1.038038146999685|1.050694694915585|1.038672188381092",
    "15|1.02|1.030706442024361|1.038881973162913|1.034570217901536|1.047266814293543|1.037079478092024|1.037275819883936|1.042422697445041|1.038126826716357|1.050777062357591|1.03874887052537",
    "16|1.02|1.031372709403557|1.03938824859539|1.03520258101596|1.047860519334206|1.037233251527569|1.037721353397281|1.042813728978499|1.038642937980465|1.051256288453125|1.039195036747409",
    "17|1.02|1.031790681397894|1.039705762654914|1.035599405964631|1.048232973179104|1.037329147928975|1.038000686627323|1.043058738985453|1.038966642954741|1.05155672843932|1.039474766662676",
    "18|1.02|1.032034491963892|1.03989094175017|1.035830924954154|1.048450232885652|1.03738488020182|1.038163567142466|1.04320155042724|1.039155440858581|1.051731910103955|1.039637878486835",
    "19|1.02|1.032117627550694|1.03995407937612|1.035909876730612|1.048524315137871|1.03740384911576|1.038219096733056|1.043250228679391|1.039219813731874|1.051791632354174|1.039693486935812",
    "20|1.02|1.031745835454002|1.039671698619523|1.035556824432633|1.048193010988073|1.037318880093845|1.037970721968757|1.043032461942372|1.038931913893744|1.051524500264818|1.039444759450859",
    "21|1.02|1.030537273755629|1.038753400159971|1.03440969591426|1.04711607095865|1.037040258392919|1.03716264616615|1.042323321271331|1.037995762590565|1.050655324841992|1.038635536087924",
    "22|1.02|1.029777531497692|1.038175839535877|1.033688971657596|1.04643908232678|1.036863255452279|1.036654128682842|1.04187656881567|1.037407044287026|1.050108302595477|1.038126296451462",
    "23|1.02|1.030180273205584|1.038482033180093|1.034070991789539|1.046797954689316|1.03695726072352|1.036923745864578|1.04211348477722|1.037719146598713|1.050398340237324|1.038396296520522",
    "24|1.02|1.031766099360422|1.039687090757186|1.035576065015612|1.04821106813971|1.037323520312145|1.03798426186202|1.04304433571211|1.038947606503832|1.051539062974626|1.039458318572323",
    "25|1.02|1.033606935384708|1.041084648215526|1.037324881922585|1.049851442553436|1.037740600503547|1.039212976322549|1.044120661446663|1.040372976732549|1.052860457201773|1.040688777948321"
)

foreach ($line in $data) {
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $val = [double]$parts[$i + 1]
        $ws.Range($cols[$i] + $row).Value = $val
    }
}
